$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Expected Launch : Mar 2026"
$ws.Range("C3").Value = "Expected Launch : Mar 2026"
